# Applies the "Add files via upload" commit:
#   - Removes the "Kai-Tak" and "Nock-Ten" rows from the "Extra names" sheet
#     (these names are dropped from the extra-names list, and their now-
#     unused shared-string entries are pruned automatically on save).
#   - Updates the description of the "Extra Names" list on the
#     "Information" sheet.
#   - Leaves the workbook with "Extra names" as the active sheet/tab, with
#     cell F19 selected there, and cell A13 selected on "Information".

$wb = $excel.ActiveWorkbook

$wsExtra = $wb.Worksheets.Item("Extra names")
$wsInfo  = $wb.Worksheets.Item("Information")

# Row 52 holds "Kai-Tak"; row 79 holds "Nock-Ten" (before any shifting).
# Deleting row 52 first shifts "Nock-Ten" up to row 78.
$wsExtra.Rows(52).Delete()
$wsExtra.Rows(78).Delete()

# Update the "Extra Names" description text (column B, row 13).
$wsInfo.Range("B13").Value = "Extra names. Initially constructed as the intersection of names in the file ``Extra-names.xlsx`` and those in the file ``IBTRACS-names.xlsx``. To be augmented according to user feedback."

# Restore view/selection state: Information!A13 selected, then make
# "Extra names" the active sheet with F19 selected (so it ends up the
# workbook's active tab).
$wsInfo.Activate()
$wsInfo.Range("A13").Select()

$wsExtra.Activate()
$wsExtra.Range("F19").Select()
